# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet and
# moves the "latest row" date-only formatting down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80 was previously the last row (formatted as date-only, style s="3").
# Since it is no longer the last row, give it the regular datetime number
# format used by all other non-final rows (style s="2").
$ws.Range("A80").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new last row (81) with the next day's values, and format its
# date cell with the date-only number format reserved for the last row.
$ws.Range("A81").Value = 45668
$ws.Range("A81").NumberFormat = "YYYY-MM-DD"
$ws.Range("B81").Value = 190
$ws.Range("C81").Value = 186
$ws.Range("D81").Value = 188
